$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold, border, centered) and date formatting
# from the "Weekly Quantity" sheet so the new sheet's styles match.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A30").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = @(
  @(44955.99999999999, 122, 13.11243852860275, 229.4650726286971),
  @(44976.99999999999, 121, 15.2045297926248, 225.6619904694245),
  @(44983.99999999999, 121, 7.790564752204543, 224.9422920160023),
  @(45039.99999999999, 118, 12.86600906581908, 229.292316787188),
  @(45046.99999999999, 118, 1.120963724037024, 229.4820235704458),
  @(45053.99999999999, 117, 4.215455363283253, 217.9954665142484),
  @(45088.99999999999, 115, 8.491661805822826, 226.7886586968094),
  @(45165.99999999999, 112, -3.561690702316346, 231.1971105944764),
  @(45179.99999999999, 111, 5.595597938840538, 214.596350831855),
  @(45186.99999999999, 110, 7.701594686501076, 215.6111456907157),
  @(45193.99999999999, 110, -7.199893217164924, 213.8476485385872),
  @(45200.99999999999, 110, -8.730762380811663, 219.6559504328542),
  @(45207.99999999999, 109, 6.891982046012372, 221.1982895646657),
  @(45214.99999999999, 109, -0.1565462564405934, 216.6898075642241),
  @(45221.99999999999, 109, -5.63643055573695, 224.0780976301795),
  @(45228.99999999999, 108, 4.760066011427082, 210.2308106258598),
  @(45235.99999999999, 108, -3.705113346036154, 217.0949077404881),
  @(45242.99999999999, 108, 0.7512575800990327, 217.0984037408506),
  @(45249.99999999999, 107, 10.4758388235612, 216.2543948204128),
  @(45256.99999999999, 107, -7.681180061340684, 222.1296214534354),
  @(45263.99999999999, 107, 0.8287951077818696, 216.9440831094501),
  @(45270.99999999999, 106, -4.316437516135844, 221.7040087438546),
  @(45277.99999999999, 106, 1.12368816476789, 213.6778438484169),
  @(45284.99999999999, 105, 2.071023197233975, 211.7807343407035),
  @(45291.99999999999, 105, 2.304706151336144, 207.9684287177844),
  @(45298.99999999999, 105, -7.500730384712461, 217.69307295035),
  @(45305.99999999999, 104, -4.74284894436361, 215.6096101814055),
  @(45312.99999999999, 104, 0.4956776041655725, 208.2934525040605),
  @(45319.99999999999, 104, -4.356741251384051, 210.970187841945)
)

$row = 2
foreach ($r in $data) {
  $wsForecast.Cells.Item($row, 1).Value = $r[0]
  $wsForecast.Cells.Item($row, 2).Value = $r[1]
  $wsForecast.Cells.Item($row, 3).Value = $r[2]
  $wsForecast.Cells.Item($row, 4).Value = $r[3]
  $row = $row + 1
}

# Restore original active sheet/selection so the workbook-level view state
# (activeTab) is unaffected by the edits made above.
$wsWeekly.Activate()
[void]$wsWeekly.Range("A1").Select()
